$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text representation
# (values like "26.619.91" or "0.9990" must stay literal text, not become numbers)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.619.91'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.53'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.97'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -6.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9992'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5272'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3303'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06760'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.52'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7794'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07666'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.836.75'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.10'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.084'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9982'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.20'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007932'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.636.08'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.073.08'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.612'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.679'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.022'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.358'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '144.81'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.647'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.05'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '111.44'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.250'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.210'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08797'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04870'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.148'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.856'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7130'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.120'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01815'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.250'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4988'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '113.49'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9060'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.098'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.862'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9991'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4319'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1296'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.157'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05928'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.44'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.441'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.88%  '
